$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "demand2" and "net2" rows/columns from the 7x7 matrix,
# shrinking it to a 5x5 demand1/net1/pv1/bat1 matrix. Delete from the
# highest index down so earlier deletions don't shift later targets.
$ws.Range("E1").EntireColumn.Delete()
$ws.Range("C1").EntireColumn.Delete()
$ws.Range("A5").EntireRow.Delete()
$ws.Range("A3").EntireRow.Delete()
